$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1127.2646
$ws.Range("I98").Value = 1025.0605
$ws.Range("J98").Value = 4500
$ws.Range("K98").Value = 1025.0605
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = 472.9395
$ws.Range("N98").Value = -7496

$ws.Range("H99").Value = 1667722.6
$ws.Range("I99").Value = 343.83334
$ws.Range("K99").Value = 1031.50002
$ws.Range("M99").Value = 466.4999800000001

$ws.Range("H106").Value = 2957.6785
$ws.Range("I106").Value = 2304.5
$ws.Range("K106").Value = 2304.5
$ws.Range("M106").Value = -1673.5

$ws.Range("H122").Value = 1127.2646
$ws.Range("I122").Value = 1025.0605
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3075.1815
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -625.1815000000001
$ws.Range("N122").Value = -18400

$ws.Range("H124").Value = 59244.57
$ws.Range("J124").Value = 59244.57
$ws.Range("L124").Value = 59244.57
$ws.Range("N124").Value = -69064.57000000001

$ws.Range("H135").Value = 200001010
$ws.Range("I135").Value = 1266
$ws.Range("K135").Value = 11394
$ws.Range("M135").Value = -8859

$ws.Range("H138").Value = 3278.2642
$ws.Range("J138").Value = 3190.7856
$ws.Range("L138").Value = 9572.356800000001
$ws.Range("N138").Value = -19852.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 397831.34
$ws.Range("I2").Value = 926614.7
$ws.Range("K2").Value = 926614.7
$ws.Range("M2").Value = -926501.7

$ws.Range("H32").Value = 14303.842
$ws.Range("I32").Value = 11986.219
$ws.Range("K32").Value = 11986.219
$ws.Range("M32").Value = -11699.219

$ws.Range("H116").Value = 397831.34
$ws.Range("I116").Value = 926614.7
$ws.Range("K116").Value = 926614.7
$ws.Range("M116").Value = -924320.7

$ws.Range("H132").Value = 1863.1945
$ws.Range("I132").Value = 1489.2046
$ws.Range("K132").Value = 4467.6138
$ws.Range("M132").Value = -1937.6138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 397831.34
$ws.Range("I3").Value = 926614.7
$ws.Range("K3").Value = 926614.7
$ws.Range("M3").Value = -926500.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3389.9167
$ws.Range("I62").Value = 2864.3333
$ws.Range("J62").Value = 4966.6665
$ws.Range("K62").Value = 2864.3333
$ws.Range("L62").Value = 4966.6665
$ws.Range("M62").Value = -2240.3333
$ws.Range("N62").Value = -6214.6665

$ws.Range("H65").Value = 3389.9167
$ws.Range("I65").Value = 2864.3333
$ws.Range("J65").Value = 4966.6665
$ws.Range("K65").Value = 14321.6665
$ws.Range("L65").Value = 24833.3325
$ws.Range("M65").Value = -11201.6665
$ws.Range("N65").Value = -31073.3325

$ws.Range("H74").Value = 34998.832
$ws.Range("J74").Value = 34998.832
$ws.Range("L74").Value = 34998.832
$ws.Range("N74").Value = -36746.832

$ws.Range("H77").Value = 34998.832
$ws.Range("J77").Value = 34998.832
$ws.Range("L77").Value = 104996.496
$ws.Range("N77").Value = -113732.496

$ws.Range("H107").Value = 1342.3448
$ws.Range("I107").Value = 1306.3182
$ws.Range("J107").Value = 1455.5714
$ws.Range("K107").Value = 1306.3182
$ws.Range("L107").Value = 1455.5714
$ws.Range("M107").Value = 613.6818000000001
$ws.Range("N107").Value = -5295.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 197.5
$ws.Range("J2").Value = 218.57143
$ws.Range("L2").Value = 1311.42858
$ws.Range("N2").Value = -1537.42858

$ws.Range("H97").Value = 919.9167
$ws.Range("I97").Value = 248.77777
$ws.Range("K97").Value = 746.33331
$ws.Range("M97").Value = -250.33331

$ws.Range("H113").Value = 47038.5
$ws.Range("J113").Value = 1305.4117
$ws.Range("L113").Value = 3916.2351
$ws.Range("N113").Value = -8256.2351

$ws.Range("H122").Value = 1392.4615
$ws.Range("I122").Value = 519.8
$ws.Range("K122").Value = 4678.2
$ws.Range("M122").Value = -2228.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1037.3334
$ws.Range("I113").Value = 590.6667
$ws.Range("K113").Value = 590.6667
$ws.Range("M113").Value = 1579.3333

$ws.Range("H132").Value = 804698.8
$ws.Range("I132").Value = 919228.1
$ws.Range("J132").Value = 2993.6667
$ws.Range("K132").Value = 2757684.3
$ws.Range("L132").Value = 8981.000100000001
$ws.Range("M132").Value = -2755154.3
$ws.Range("N132").Value = -14041.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3020.8667
$ws.Range("I7").Value = 2766.889
$ws.Range("J7").Value = 3401.8333
$ws.Range("K7").Value = 2766.889
$ws.Range("L7").Value = 3401.8333
$ws.Range("M7").Value = -2654.889
$ws.Range("N7").Value = -3625.8333

$ws.Range("H40").Value = 14500.125
$ws.Range("I40").Value = 14789.223
$ws.Range("J40").Value = 14128.429
$ws.Range("K40").Value = 14789.223
$ws.Range("L40").Value = 14128.429
$ws.Range("M40").Value = -14653.223
$ws.Range("N40").Value = -14400.429

$ws.Range("H46").Value = 2045.3572
$ws.Range("J46").Value = 2268.25
$ws.Range("L46").Value = 2268.25
$ws.Range("N46").Value = -2644.25

$ws.Range("H61").Value = 2666.1
$ws.Range("I61").Value = 2254.4119
$ws.Range("K61").Value = 2254.4119
$ws.Range("M61").Value = -2052.4119

$ws.Range("H82").Value = 1912.2727
$ws.Range("I82").Value = 1974.6666
$ws.Range("K82").Value = 1974.6666
$ws.Range("M82").Value = -1613.6666

$ws.Range("H85").Value = 1912.2727
$ws.Range("I85").Value = 1974.6666
$ws.Range("K85").Value = 1974.6666
$ws.Range("M85").Value = -726.6666

$ws.Range("H96").Value = 74993
$ws.Range("J96").Value = 74993
$ws.Range("L96").Value = 74993
$ws.Range("N96").Value = -80485

$ws.Range("H113").Value = 2666.1
$ws.Range("I113").Value = 2254.4119
$ws.Range("K113").Value = 2254.4119
$ws.Range("M113").Value = -84.41190000000006

$ws.Range("H126").Value = 3020.8667
$ws.Range("I126").Value = 2766.889
$ws.Range("J126").Value = 3401.8333
$ws.Range("K126").Value = 8300.667000000001
$ws.Range("L126").Value = 10205.4999
$ws.Range("M126").Value = -5830.667000000001
$ws.Range("N126").Value = -15145.4999

$ws.Range("H132").Value = 3606.8333
$ws.Range("I132").Value = 3308.5454
$ws.Range("J132").Value = 3934.95
$ws.Range("K132").Value = 9925.636200000001
$ws.Range("L132").Value = 11804.85
$ws.Range("M132").Value = -7395.636200000001
$ws.Range("N132").Value = -16864.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 919.7646999999999
$ws.Range("I113").Value = 820.5
$ws.Range("J113").Value = 1061.5714
$ws.Range("K113").Value = 2461.5
$ws.Range("L113").Value = 3184.7142
$ws.Range("M113").Value = -291.5
$ws.Range("N113").Value = -7524.7142

$ws.Range("H126").Value = 3420.4
$ws.Range("I126").Value = 3420.4
$ws.Range("K126").Value = 10261.2
$ws.Range("M126").Value = -7791.200000000001

$ws.Range("H132").Value = 1093.9231
$ws.Range("I132").Value = 887.26666
$ws.Range("K132").Value = 2661.79998
$ws.Range("M132").Value = -131.7999799999998
